$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: bulk-upload fix drops the "createdBy" column (old F), shifting
# "location" (old G) into column F.
$data = @(
    @("name", "email", "role", "level", "progress", "location"),
    @("Joshua Black", "mhart@example.com", "UI/UX Designer", "MID", "OFFER ACCEPTED", "Watsontown"),
    @("Daniel Cunningham", "derekbell@example.net", "Backend Developer", "MID", "OFFER REJECTED", "Curtisstad"),
    @("Tammy David", "robinsonsara@example.org", "DevOps Engineer", "JUNIOR", "SHORTLISTED", "North Keithville"),
    @("Lindsey Williams", "cobbwilliam@example.com", "Backend Developer", "MID", "REJECTED", "Milesland"),
    @("Paul Reed", "deborahlucas@example.com", "DevOps Engineer", "MID", "OFFERED", "South Johnton"),
    @("Jamie Liu", "ggallagher@example.org", "Product Manager", "SENIOR", "ON HOLD", "Sanchezside"),
    @("Nancy Mann", "ashley35@example.org", "Frontend Developer", "SENIOR", "REJECTED", "North Brittneyshire"),
    @("Crystal Kaufman", "sspencer@example.org", "Backend Developer", "LEAD", "REJECTED", "Khanport"),
    @("Amanda Velazquez", "bethsmith@example.com", "DevOps Engineer", "LEAD", "OFFER REJECTED", "Johnland"),
    @("Paul Owens", "melissachavez@example.net", "Backend Developer", "SENIOR", "REJECTED", "Zacharyfort")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Old column G (createdBy header / data) is no longer part of the export;
# drop it so the sheet shrinks back to 6 columns.
$ws.Columns("G").Delete()
